$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3ASY04_Genomics")

# --- Header (building block type) renames in row 1 / table column names ---
$ws.Range("B1").Value = "Characteristic [BioSample Accession Number]"
$ws.Range("H1").Value = "Characteristic [library source]"
$ws.Range("Q1").Value = "Component [library preparation kit]"
$ws.Range("AG1").Value = "Component [next generation sequencing instrument model]"
$ws.Range("AJ1").Value = "Component [base-calling software]"

# --- Data row (row 2) value updates ---
# "library source" example value: genomic DNA -> Genomic DNA
$ws.Range("H2").Value = "Genomic DNA"

# Term Source REF / Term Accession Number for the library source value now
# point at NCIT:C95940 instead of BAO:0000316
$ws.Range("I2").Value = "NCIT"
$ws.Range("J2").Value = "https://bioregistry.io/NCIT:C95940"

# Term Source REF / Term Accession Number for the sequencing instrument model
# platform now point at EFO:0004203 instead of OBI:0002001
$ws.Range("AH2").Value = "EFO"
$ws.Range("AI2").Value = "https://bioregistry.io/EFO:0004203"
